$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 91, shifting existing rows 91:204 down to 92:205.
$ws.Rows(91).Insert()

# Populate the newly inserted row 91 with the new daily record.
$ws.Range("A91").Value = 5
$ws.Range("B91").Value = "Macroferia Regional de Talca"
$ws.Range("C91").Value = "Maule"
$ws.Range("D91").Value = 44483
$ws.Range("E91").Value = 7
$ws.Range("F91").Value = 100112023
$ws.Range("G91").Value = "Brócoli"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 4000
$ws.Range("K91").Value = 800
$ws.Range("L91").Value = 800
$ws.Range("M91").Value = 800
$ws.Range("N91").Value = "$/unidad"
$ws.Range("O91").Value = "Región Metropolitana"
$ws.Range("P91").Value = 800
$ws.Range("Q91").Value = 1
$ws.Range("R91").Value = "Hortaliza"
